$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.215.57"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "2.323.41"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'509.33"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'131.96"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "'0.0995"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'5.24"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "2.738.23"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'23.38"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "56.202.08"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "2.328.68"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'10.34"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "'321.47"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'61.12"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'8.54"
$ws.Range("E24").Value = "  +10.12%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.161"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'1.29"
$ws.Range("E27").Value = "  +4.25%  "
$ws.Range("D28").Value = "'166.97"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "'1.66"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "0.0₃0712"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "'0.879"
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("D38").Value = "'38.40"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.53"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'149.69"
$ws.Range("E40").Value = "  +9.11%  "
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "'3.54"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "'276.62"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'4.98"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "'0.0919"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'0.0493"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "'0.0213"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'16.92"
$ws.Range("E51").Value = "  +0.45%  "
